$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.723979
$ws.Range("H2").Value = 23.171937
$ws.Range("I2").Value = 0.471042132528101
$ws.Range("J2").Value = 0.471042132528101
$ws.Range("M2").Value = 0.3951603333333333
$ws.Range("N2").Value = 1.185481
$ws.Range("O2").Value = 0.02514070644417849
$ws.Range("P2").Value = 0.02514070644417849
$ws.Range("Q2").Value = 3.052210116299667
$ws.Range("R2").Value = 27.469891046697
$ws.Range("S2").Value = 0.0118423319767288
$ws.Range("T2").Value = 0.0118423319767288

$ws.Range("G3").Value = 7.723979
$ws.Range("H3").Value = 23.171937
$ws.Range("I3").Value = 0.471042132528101
$ws.Range("J3").Value = 0.471042132528101
$ws.Range("O3").Value = 0.1086999987233279
$ws.Range("P3").Value = 0.1086999987233279
$ws.Range("Q3").Value = 13.196734804639
$ws.Range("R3").Value = 118.770613241751
$ws.Range("S3").Value = 0.05120227920443823
$ws.Range("T3").Value = 0.05120227920443822

$ws.Range("G4").Value = 7.723979
$ws.Range("H4").Value = 23.171937
$ws.Range("I4").Value = 0.471042132528101
$ws.Range("J4").Value = 0.471042132528101
$ws.Range("O4").Value = 0.8661592948324937
$ws.Range("P4").Value = 0.8661592948324937
$ws.Range("Q4").Value = 105.1561605034727
$ws.Range("R4").Value = 946.405444531254
$ws.Range("S4").Value = 0.407997521346934
$ws.Range("T4").Value = 0.407997521346934

$ws.Range("I5").Value = 0.2460132574367717
$ws.Range("J5").Value = 0.2460132574367717
$ws.Range("M5").Value = 0.3951603333333333
$ws.Range("N5").Value = 1.185481
$ws.Range("O5").Value = 0.02514070644417849
$ws.Range("P5").Value = 0.02514070644417849
$ws.Range("Q5").Value = 1.594091273878889
$ws.Range("R5").Value = 14.34682146491
$ws.Range("S5").Value = 0.006184947086593988
$ws.Range("T5").Value = 0.006184947086593988

$ws.Range("I6").Value = 0.2460132574367717
$ws.Range("J6").Value = 0.2460132574367717
$ws.Range("O6").Value = 0.1086999987233279
$ws.Range("P6").Value = 0.1086999987233279
$ws.Range("S6").Value = 0.02674164076929882
$ws.Range("T6").Value = 0.02674164076929882

$ws.Range("I7").Value = 0.2460132574367717
$ws.Range("J7").Value = 0.2460132574367717
$ws.Range("O7").Value = 0.8661592948324937
$ws.Range("P7").Value = 0.8661592948324937
$ws.Range("S7").Value = 0.2130866695808789
$ws.Range("T7").Value = 0.2130866695808789

$ws.Range("I8").Value = 0.2829446100351274
$ws.Range("J8").Value = 0.2829446100351274
$ws.Range("M8").Value = 0.3951603333333333
$ws.Range("N8").Value = 1.185481
$ws.Range("O8").Value = 0.02514070644417849
$ws.Range("P8").Value = 0.02514070644417849
$ws.Range("Q8").Value = 1.833395234661222
$ws.Range("R8").Value = 16.500557111951
$ws.Range("S8").Value = 0.007113427380855698
$ws.Range("T8").Value = 0.007113427380855698

$ws.Range("I9").Value = 0.2829446100351274
$ws.Range("J9").Value = 0.2829446100351274
$ws.Range("O9").Value = 0.1086999987233279
$ws.Range("P9").Value = 0.1086999987233279
$ws.Range("S9").Value = 0.03075607874959086
$ws.Range("T9").Value = 0.03075607874959086

$ws.Range("I10").Value = 0.2829446100351274
$ws.Range("J10").Value = 0.2829446100351274
$ws.Range("O10").Value = 0.8661592948324937
$ws.Range("P10").Value = 0.8661592948324937
$ws.Range("S10").Value = 0.2450751039046809
$ws.Range("T10").Value = 0.2450751039046809
